$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$tbl = $ws.ListObjects.Item("Table1")

Write-Host ("Before: " + $tbl.Range.Address())

$ws.Rows.Item(87).Insert()
$ws.Rows.Item(89).Insert()
$ws.Rows.Item(91).Insert()
$ws.Rows.Item(92).Insert()
$ws.Rows.Item(93).Insert()

Write-Host ("After inserts (sheet rows), table range unchanged?: " + $tbl.Range.Address())

$tbl.Resize($ws.Range("A8:K121"))
Write-Host ("After resize: " + $tbl.Range.Address())

for ($r = 84; $r -le 94; $r++) {
    Write-Host ("Row " + $r + " A=" + $ws.Cells.Item($r,1).Value2 + " style A=" + $ws.Cells.Item($r,1).Value2)
}
